# Daily attendance processing - 2026-02-01 11:11:31
# Normalizes the "Recorded By" (column G) entries so that the "System"
# author is listed first in the comma-separated author list, instead of
# last. Entries that don't contain an exact "System" token, or that were
# recorded by an admin account, are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value()

    if ($null -eq $val -or $val -eq "") {
        continue
    }

    # Admin-recorded rows are left alone.
    if ($val.Contains("admin@admin.com")) {
        continue
    }

    $parts = $val -split ", "

    # Find the last exact "System" token (case-sensitive) in the list.
    $lastSystemIndex = -1
    for ($i = 0; $i -lt $parts.Length; $i++) {
        if ($parts[$i] -ceq "System") {
            $lastSystemIndex = $i
        }
    }

    if ($lastSystemIndex -lt 0) {
        continue
    }

    $newParts = @($parts[$lastSystemIndex])
    for ($i = 0; $i -lt $parts.Length; $i++) {
        if ($i -ne $lastSystemIndex) {
            $newParts += $parts[$i]
        }
    }

    $newVal = [string]::Join(", ", $newParts)

    if ($newVal -cne $val) {
        $cell.Value = $newVal
    }
}
